$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.018.60"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "1.867.25"
$ws.Range("E3").Value = "  -1.13%  "
$c = $ws.Range("D4")
$c.Value = "'1.004"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.18%  "
$c = $ws.Range("D5")
$c.Value = "'312.76"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.35%  "
$c = $ws.Range("D6")
$c.Value = "'1.004"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.24%  "
$c = $ws.Range("D7")
$c.Value = "'0.5078"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +1.54%  "
$c = $ws.Range("D8")
$c.Value = "'0.3806"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.79%  "
$c = $ws.Range("D9")
$c.Value = "'0.08316"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -9.28%  "
$c = $ws.Range("D10")
$c.Value = "'1.107"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.68%  "
$c = $ws.Range("D11")
$c.Value = "'41.38"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.82%  "
$c = $ws.Range("D12")
$c.Value = "'6.206"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -2.38%  "
$ws.Range("D13").Value = "1.868.98"
$ws.Range("E13").Value = "  -0.95%  "
$c = $ws.Range("D14")
$c.Value = "'20.45"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.52%  "
$c = $ws.Range("D15")
$c.Value = "'7.181"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.58%  "
$ws.Range("E16").Value = "  +0.21%  "
$c = $ws.Range("D17")
$c.Value = "'0.00001094"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.95%  "
$c = $ws.Range("D18")
$c.Value = "'90.48"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.27%  "
$c = $ws.Range("D19")
$c.Value = "'0.06639"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.13%  "
$c = $ws.Range("D20")
$c.Value = "'17.84"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("E21").Value = "  +0.14%  "
$c = $ws.Range("D22")
$c.Value = "'6.010"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -3.48%  "
$ws.Range("D23").Value = "28.045.39"
$ws.Range("E23").Value = "  -0.29%  "
$c = $ws.Range("D24")
$c.Value = "'11.13"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.92%  "
$c = $ws.Range("D25")
$c.Value = "'2.259"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -2.48%  "
$c = $ws.Range("D26")
$c.Value = "'2.569"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.83%  "
$ws.Range("D27").Value = "2.080.12"
$ws.Range("E27").Value = "  -1.14%  "
$c = $ws.Range("D28")
$c.Value = "'157.49"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.58%  "
$c = $ws.Range("D29")
$c.Value = "'20.51"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.07%  "
$c = $ws.Range("D30")
$c.Value = "'126.25"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.39%  "
$c = $ws.Range("D31")
$c.Value = "'0.1056"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.36%  "
$c = $ws.Range("D32")
$c.Value = "'1.038"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -2.95%  "
$c = $ws.Range("D33")
$c.Value = "'5.583"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.05%  "
$c = $ws.Range("D34")
$c.Value = "'3.596"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.02%  "
$c = $ws.Range("D35")
$c.Value = "'9.635"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +2.41%  "
$c = $ws.Range("D36")
$c.Value = "'0.02424"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.70%  "
$ws.Range("E37").Value = "  -0.83%  "
$ws.Range("E38").Value = "  -1.66%  "
$c = $ws.Range("D39")
$c.Value = "'1.207"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.46%  "
$c = $ws.Range("D40")
$c.Value = "'0.6397"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("E41").Value = "  -6.55%  "
$c = $ws.Range("D42")
$c.Value = "'11.26"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -2.46%  "
$c = $ws.Range("D43")
$c.Value = "'4.854"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -1.73%  "
$c = $ws.Range("D44")
$c.Value = "'0.6082"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.91%  "
$c = $ws.Range("D45")
$c.Value = "'13.03"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -2.70%  "
$c = $ws.Range("D46")
$c.Value = "'1.288"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -1.00%  "
$c = $ws.Range("D47")
$c.Value = "'3.646"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -1.01%  "
$ws.Range("E48").Value = "  -0.20%  "
$c = $ws.Range("D49")
$c.Value = "'1.209"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.32%  "
$c = $ws.Range("D50")
$c.Value = "'121.12"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.00%  "
$c = $ws.Range("D51")
$c.Value = "'79.68"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.95%  "
